$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 79; this shifts the existing rows 79-177
# down to 80-178 (and the sheet's used range grows to A1:T178).
$ws.Rows("79:79").Insert()

# Populate the newly inserted row 79 with the new weekly data point.
$ws.Range("A79").Value = 5
$ws.Range("B79").Value = "Macroferia Regional de Talca"
$ws.Range("C79").Value = "Maule"
$ws.Range("D79").Value = 45117
$ws.Range("E79").Value = 7
$ws.Range("F79").Value = "Fruta"
$ws.Range("G79").Value = 100108
$ws.Range("H79").Value = "Tropicales y subtropicales"
$ws.Range("I79").Value = 100108002
$ws.Range("J79").Value = "Mango"
$ws.Range("K79").Value = "Sin especificar"
$ws.Range("L79").Value = "Primera"
$ws.Range("M79").Value = 248
$ws.Range("N79").Value = 8000
$ws.Range("O79").Value = 8000
$ws.Range("P79").Value = 8000
$ws.Range("Q79").Value = "$/bandeja 4 kilos"
$ws.Range("R79").Value = "Brasil"
$ws.Range("S79").Value = 2000
$ws.Range("T79").Value = 4
